# Update the "correlations sp500" worksheet:
#  - Insert a new "bank" keyword row (with its S&P 500 correlation) right
#    after the header row, pushing credit/inflation/interest down by one row.
#  - Drop the "mortgage" and "recession" keyword rows entirely.
#  - "trade" keeps its correlation value but now lands on row 6 instead of 7.
#
# Net effect: table shrinks from A1:C7 to A1:C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the new "bank" entry (previously held "credit").
$ws.Range("A2").Value = "bank"
$ws.Range("B2").Value = "Frequency"
$ws.Range("C2").Value = -0.2146

# Rows 3-5 shift up to hold credit / inflation / interest (mortgage/recession removed).
$ws.Range("A3").Value = "credit"
$ws.Range("B3").Value = "Frequency"
$ws.Range("C3").Value = -0.2897

$ws.Range("A4").Value = "inflation"
$ws.Range("B4").Value = "Frequency"
$ws.Range("C4").Value = 0.3226

$ws.Range("A5").Value = "interest"
$ws.Range("B5").Value = "Frequency"
$ws.Range("C5").Value = -0.0175

# "trade" now occupies row 6 (mortgage/recession rows 5-6 are gone).
$ws.Range("A6").Value = "trade"
$ws.Range("B6").Value = "Frequency"
$ws.Range("C6").Value = 0.2616

# The old row 7 (trade's previous location) is no longer part of the table.
$ws.Range("A7:C7").Clear()
